# Append 12 new sign-up attempt rows (154-165) to the Users sheet.
# Row 155 is the new "Dianka" entry; the rest repeat the existing
# "moses/bro/1234/m@g.c/Male" test submission already used elsewhere
# in the sheet (e.g. row 6).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, [string]$text) {
    # Force the cell to be stored as text even when the content looks
    # like a number (e.g. "1234"), then restore the default "Normal"
    # style so no stray number-format style lingers on the cell.
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

$repeatRow = @{
    A = "moses"
    B = "bro"
    C = "1234"
    D = "m@g.c"
    E = "Male"
    F = 0
}

$dianaRow = @{
    A = "Dianka"
    B = "D12345!1"
    C = "123456789"
    D = "Diana@hey.com"
    E = "Female"
    F = 0
}

function Write-UserRow([int]$rowNum, $data) {
    Set-TextValue $ws.Cells.Item($rowNum, 1) $data.A
    Set-TextValue $ws.Cells.Item($rowNum, 2) $data.B
    Set-TextValue $ws.Cells.Item($rowNum, 3) $data.C
    Set-TextValue $ws.Cells.Item($rowNum, 4) $data.D
    Set-TextValue $ws.Cells.Item($rowNum, 5) $data.E
    $ws.Cells.Item($rowNum, 6).Value = $data.F
}

Write-UserRow 154 $repeatRow
Write-UserRow 155 $dianaRow
Write-UserRow 156 $repeatRow
Write-UserRow 157 $repeatRow
Write-UserRow 158 $repeatRow
Write-UserRow 159 $repeatRow
Write-UserRow 160 $repeatRow
Write-UserRow 161 $repeatRow
Write-UserRow 162 $repeatRow
Write-UserRow 163 $repeatRow
Write-UserRow 164 $repeatRow
Write-UserRow 165 $repeatRow
